$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the fit-statistic strings between B4/B6 and B5/B7
$b4 = $ws.Range("B4").Text
$b5 = $ws.Range("B5").Text
$b6 = $ws.Range("B6").Text
$b7 = $ws.Range("B7").Text

$ws.Range("B4").Value = $b6
$ws.Range("B5").Value = $b7
$ws.Range("B6").Value = $b4
$ws.Range("B7").Value = $b5
